$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: insert a new "Meta description" paragraph right after the
# title (Heading1) paragraph. Use InsertXML so the run layout
# (leading empty run + bold "Meta description" run + plain run with
# the rest of the sentence) matches exactly what Word itself produces.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.First
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Style = "Normal"

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore a medieval realm with Cold Spell, a Novomatic online slot game with stunning ice-covered reels and exciting win potential. Play for free now.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------
# Step 2: remove the duplicate bold "Play Cold Spell Slot for Free -
# Novomatic Fantasy Theme" paragraph that used to sit just before the
# closing (italic) paragraph near the end of the document. Locate it
# by content (skipping the real Heading1 title, paragraph 1) so the
# removal does not depend on a fixed paragraph count/index.
# ---------------------------------------------------------------------
$titleText = "Play Cold Spell Slot for Free"
$boldParaIndex = -1
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*$titleText*") {
        $boldParaIndex = $i
    }
}
$boldPara = $d.Paragraphs($boldParaIndex)
$boldPara.Range.Delete()

# ---------------------------------------------------------------------
# Step 3: rewrite the closing italic paragraph's text (keeps the
# existing italic run formatting) with the new image-prompt copy.
# ---------------------------------------------------------------------
$closingPara = $d.Paragraphs($d.Paragraphs.Count)
$closingRange = $d.Range($closingPara.Range.Start, $closingPara.Range.End)
$closingRange.Text = "Please create a cartoon-style feature image for the online slot game " + [char]34 + "Cold Spell" + [char]34 + ". The image should feature a happy Maya warrior with glasses. The Maya warrior should be depicted holding a wand and standing in front of ice-covered mountains with a snowy background to reflect the game" + [char]39 + "s medieval fantasy theme. The image can include other elements from the game such as playing cards, tiaras, maps, and treasure chests. The image should be bright and colorful to capture the attention of players and entice them to try the game."
